# Border_Flows.xlsx update
# - Shift every timestamp in column A (rows 2-101) forward by exactly 1 day
#   (the ENTSO-E fetch window rolled forward one day).
# - Zero out the flow values in B6:N33 (the "PC Sun model" forecast rows
#   that produced bad data and need retraining, per the commit message).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 101

# Shift column A (Timestamp) down by one full day for every data row.
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value2 = $cell.Value2 + 1
}

# Zero out the affected flow columns (B:N) for rows 6 through 33.
$ws.Range("B6:N33").Value2 = 0
